{"js": "// Update the 25 multiplication problems (three-digit number x one-digit\n// number) that live in the table cells of the worksheet. Each old problem\n// string is unique in the document, so we can safely match on the exact\n// text and swap in the new problem.\nconst replacements = [\n  [\"780\u00d79=\", \"372\u00d72=\"],\n  [\"873\u00d76=\", \"735\u00d76=\"],\n  [\"502\u00d79=\", \"121\u00d76=\"],\n  [\"938\u00d73=\", \"354\u00d72=\"],\n  [\"613\u00d79=\", \"854\u00d77=\"],\n  [\"834\u00d76=\", \"335\u00d78=\"],\n  [\"333\u00d79=\", \"817\u00d78=\"],\n  [\"566\u00d79=\", \"313\u00d72=\"],\n  [\"936\u00d76=\", \"627\u00d72=\"],\n  [\"439\u00d75=\", \"804\u00d76=\"],\n  [\"962\u00d78=\", \"296\u00d79=\"],\n  [\"904\u00d78=\", \"814\u00d79=\"],\n  [\"840\u00d77=\", \"175\u00d79=\"],\n  [\"995\u00d76=\", \"285\u00d74=\"],\n  [\"505\u00d76=\", \"226\u00d72=\"],\n  [\"962\u00d79=\", \"944\u00d72=\"],\n  [\"746\u00d73=\", \"933\u00d77=\"],\n  [\"877\u00d77=\", \"128\u00d78=\"],\n  [\"556\u00d74=\", \"922\u00d79=\"],\n  [\"475\u00d78=\", \"697\u00d72=\"],\n  [\"211\u00d73=\", \"516\u00d75=\"],\n  [\"618\u00d74=\", \"442\u00d73=\"],\n  [\"756\u00d76=\", \"319\u00d74=\"],\n  [\"853\u00d73=\", \"610\u00d75=\"],\n  [\"816\u00d73=\", \"865\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 multiplication problems (three-digit number x one-digit\n# number) that live in the table cells of the worksheet. Each old problem\n# string is unique in the document, so Find/Replace on the exact text is\n# unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"780\u00d79=\"; New = \"372\u00d72=\" },\n  @{ Old = \"873\u00d76=\"; New = \"735\u00d76=\" },\n  @{ Old = \"502\u00d79=\"; New = \"121\u00d76=\" },\n  @{ Old = \"938\u00d73=\"; New = \"354\u00d72=\" },\n  @{ Old = \"613\u00d79=\"; New = \"854\u00d77=\" },\n  @{ Old = \"834\u00d76=\"; New = \"335\u00d78=\" },\n  @{ Old = \"333\u00d79=\"; New = \"817\u00d78=\" },\n  @{ Old = \"566\u00d79=\"; New = \"313\u00d72=\" },\n  @{ Old = \"936\u00d76=\"; New = \"627\u00d72=\" },\n  @{ Old = \"439\u00d75=\"; New = \"804\u00d76=\" },\n  @{ Old = \"962\u00d78=\"; New = \"296\u00d79=\" },\n  @{ Old = \"904\u00d78=\"; New = \"814\u00d79=\" },\n  @{ Old = \"840\u00d77=\"; New = \"175\u00d79=\" },\n  @{ Old = \"995\u00d76=\"; New = \"285\u00d74=\" },\n  @{ Old = \"505\u00d76=\"; New = \"226\u00d72=\" },\n  @{ Old = \"962\u00d79=\"; New = \"944\u00d72=\" },\n  @{ Old = \"746\u00d73=\"; New = \"933\u00d77=\" },\n  @{ Old = \"877\u00d77=\"; New = \"128\u00d78=\" },\n  @{ Old = \"556\u00d74=\"; New = \"922\u00d79=\" },\n  @{ Old = \"475\u00d78=\"; New = \"697\u00d72=\" },\n  @{ Old = \"211\u00d73=\"; New = \"516\u00d75=\" },\n  @{ Old = \"618\u00d74=\"; New = \"442\u00d73=\" },\n  @{ Old = \"756\u00d76=\"; New = \"319\u00d74=\" },\n  @{ Old = \"853\u00d73=\"; New = \"610\u00d75=\" },\n  @{ Old = \"816\u00d73=\"; New = \"865\u00d73=\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $r.New, 2) | Out-Null\n}\n"}
